$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1587.1428
$ws.Range("I28").Value = 236.4
$ws.Range("K28").Value = 236.4
$ws.Range("M28").Value = 248.6

$ws.Range("H40").Value = 1858.3334
$ws.Range("I40").Value = 1900
$ws.Range("J40").Value = 1844.4445
$ws.Range("K40").Value = 1900
$ws.Range("L40").Value = 1844.4445
$ws.Range("M40").Value = -1725
$ws.Range("N40").Value = -2194.4445

$ws.Range("H41").Value = 4598.4287
$ws.Range("I41").Value = 90
$ws.Range("J41").Value = 6401.8
$ws.Range("K41").Value = 90
$ws.Range("L41").Value = 6401.8
$ws.Range("M41").Value = 350
$ws.Range("N41").Value = -7281.8

$ws.Range("H92").Value = 1193.65
$ws.Range("I92").Value = 1345.6666
$ws.Range("J92").Value = 965.625
$ws.Range("K92").Value = 1345.6666
$ws.Range("L92").Value = 965.625
$ws.Range("M92").Value = -97.66660000000002
$ws.Range("N92").Value = -3461.625

$ws.Range("H100").Value = 1071.1111
$ws.Range("I100").Value = 852
$ws.Range("J100").Value = 2166.6667
$ws.Range("K100").Value = 852
$ws.Range("L100").Value = 2166.6667
$ws.Range("M100").Value = -311
$ws.Range("N100").Value = -3248.6667

$ws.Range("H107").Value = 259.4
$ws.Range("I107").Value = 264.92856
$ws.Range("J107").Value = 246.5
$ws.Range("K107").Value = 264.92856
$ws.Range("L107").Value = 246.5
$ws.Range("M107").Value = 1655.07144
$ws.Range("N107").Value = -4086.5

$ws.Range("H132").Value = 3862873
$ws.Range("I132").Value = 4203619.5
$ws.Range("J132").Value = 1080
$ws.Range("K132").Value = 12610858.5
$ws.Range("L132").Value = 3240
$ws.Range("M132").Value = -12608328.5
$ws.Range("N132").Value = -8300

$ws.Range("H137").Value = 1327.9459
$ws.Range("I137").Value = 893.7406999999999
$ws.Range("J137").Value = 2500.3
$ws.Range("K137").Value = 2681.2221
$ws.Range("L137").Value = 7500.900000000001
$ws.Range("M137").Value = -131.2221
$ws.Range("N137").Value = -12600.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19388.344
$ws.Range("I32").Value = 20518.666
$ws.Range("J32").Value = 10184.286
$ws.Range("K32").Value = 20518.666
$ws.Range("L32").Value = 10184.286
$ws.Range("M32").Value = -20231.666
$ws.Range("N32").Value = -10758.286

$ws.Range("H76").Value = 40000
$ws.Range("J76").Value = 40000
$ws.Range("L76").Value = 40000
$ws.Range("N76").Value = -40676

$ws.Range("H79").Value = 40000
$ws.Range("J79").Value = 40000
$ws.Range("L79").Value = 40000
$ws.Range("N79").Value = -42340

$ws.Range("H92").Value = 48433.332
$ws.Range("J92").Value = 48433.332
$ws.Range("L92").Value = 48433.332
$ws.Range("N92").Value = -53425.332

$ws.Range("H132").Value = 3002.7827
$ws.Range("I132").Value = 3011.1
$ws.Range("J132").Value = 2947.3333
$ws.Range("K132").Value = 9033.299999999999
$ws.Range("L132").Value = 8841.999899999999
$ws.Range("M132").Value = -6503.299999999999
$ws.Range("N132").Value = -13901.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2159.8
$ws.Range("I86").Value = 1400
$ws.Range("J86").Value = 2666.3333
$ws.Range("K86").Value = 1400
$ws.Range("L86").Value = 2666.3333
$ws.Range("M86").Value = -277
$ws.Range("N86").Value = -4912.3333

$ws.Range("H89").Value = 2159.8
$ws.Range("I89").Value = 1400
$ws.Range("J89").Value = 2666.3333
$ws.Range("K89").Value = 7000
$ws.Range("L89").Value = 13331.6665
$ws.Range("M89").Value = -1384
$ws.Range("N89").Value = -24563.6665

$ws.Range("H105").Value = 2617.0688
$ws.Range("I105").Value = 2277.75
$ws.Range("K105").Value = 2277.75
$ws.Range("M105").Value = -530.75

$ws.Range("H138").Value = 70578
$ws.Range("J138").Value = 70578
$ws.Range("L138").Value = 70578
$ws.Range("N138").Value = -80858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 35300
$ws.Range("J92").Value = 35300
$ws.Range("L92").Value = 35300
$ws.Range("N92").Value = -40292

$ws.Range("H132").Value = 2377.6843
$ws.Range("I132").Value = 2228.4614
$ws.Range("K132").Value = 6685.3842
$ws.Range("M132").Value = -4155.3842

$ws.Range("H141").Value = 52727.273
$ws.Range("J141").Value = 58000
$ws.Range("L141").Value = 58000
$ws.Range("N141").Value = -68360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2799.804
$ws.Range("I131").Value = 13081.25
$ws.Range("J131").Value = 886.97675
$ws.Range("K131").Value = 39243.75
$ws.Range("L131").Value = 2660.93025
$ws.Range("M131").Value = -34203.75
$ws.Range("N131").Value = -12740.93025

$ws.Range("H140").Value = 1836.6471
$ws.Range("I140").Value = 1328.8667
$ws.Range("K140").Value = 3986.6001
$ws.Range("M140").Value = 1193.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H113").Value = 27778522
$ws.Range("I113").Value = 62500424
$ws.Range("K113").Value = 62500424
$ws.Range("M113").Value = -62498254

$ws.Range("H132").Value = 27902.846
$ws.Range("I132").Value = 49954.668
$ws.Range("J132").Value = 2175.7222
$ws.Range("K132").Value = 149864.004
$ws.Range("L132").Value = 6527.1666
$ws.Range("M132").Value = -147334.004
$ws.Range("N132").Value = -11587.1666

$ws.Range("H141").Value = 37402.9
$ws.Range("J141").Value = 37402.9
$ws.Range("L141").Value = 37402.9
$ws.Range("N141").Value = -47762.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7064.0884
$ws.Range("I132").Value = 11634.588
$ws.Range("J132").Value = 2493.5881
$ws.Range("K132").Value = 34903.764
$ws.Range("L132").Value = 7480.7643
$ws.Range("M132").Value = -32373.764
$ws.Range("N132").Value = -12540.7643

$ws.Range("H133").Value = 21891.777
$ws.Range("J133").Value = 21891.777
$ws.Range("L133").Value = 21891.777
$ws.Range("N133").Value = -26951.777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 28155.334
$ws.Range("J70").Value = 28155.334
$ws.Range("L70").Value = 28155.334
$ws.Range("N70").Value = -28785.334

$ws.Range("H73").Value = 28155.334
$ws.Range("J73").Value = 28155.334
$ws.Range("L73").Value = 28155.334
$ws.Range("N73").Value = -30339.334

$ws.Range("H132").Value = 1445.3
$ws.Range("I132").Value = 755.8333
$ws.Range("J132").Value = 2479.5
$ws.Range("K132").Value = 2267.4999
$ws.Range("L132").Value = 7438.5
$ws.Range("M132").Value = 262.5001000000002
$ws.Range("N132").Value = -12498.5

$ws.Range("H140").Value = 40832
$ws.Range("J140").Value = 40832
$ws.Range("L140").Value = 40832
$ws.Range("N140").Value = -51192

$ws.Range("H141").Value = 52907.5
$ws.Range("J141").Value = 52907.5
$ws.Range("L141").Value = 52907.5
$ws.Range("N141").Value = -63267.5
